{"js": "// Split the Title paragraph and the Abstract paragraph into one run per\n// word, with each inter-word space as its own run (matching the \"one\n// word per <w:r>\" convention already used elsewhere in this document).\n//\n// A plain Range.insertText()/paragraph.clear() rebuild would cause the\n// host to coalesce same-formatted adjacent runs back into a single run\n// on save, so we instead build the exact run-level OOXML for each\n// paragraph and swap it in with insertOoxml (Word.InsertLocation.replace),\n// which preserves run boundaries verbatim.\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// \"A B (C)\" -> one <w:r> per word, plus one <w:r> per single space\n// between words (never a leading/trailing space run).\nfunction wordRunsXml(text) {\n  const words = text.split(\" \");\n  let xml = \"\";\n  for (let i = 0; i < words.length; i++) {\n    if (i > 0) {\n      xml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>';\n    }\n    xml += '<w:r><w:t xml:space=\"preserve\">' + escapeXml(words[i]) + \"</w:t></w:r>\";\n  }\n  return xml;\n}\n\nfunction paragraphOoxml(pStyle, text) {\n  const runs = wordRunsXml(text);\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:pPr><w:pStyle w:val=\\\"\" +\n    pStyle +\n    '\"/></w:pPr>' +\n    runs +\n    \"</w:p></w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"style\");\nawait context.sync();\n\nlet titlePara = null;\nlet abstractPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (titlePara === null && p.style === \"Title\") {\n    titlePara = p;\n  }\n  if (abstractPara === null && p.style === \"Abstract\") {\n    abstractPara = p;\n  }\n}\n\nif (titlePara !== null) {\n  titlePara.insertOoxml(\n    paragraphOoxml(\"Title\", \"Answers: Trigonometry (degrees)\"),\n    Word.InsertLocation.replace\n  );\n}\n\nif (abstractPara !== null) {\n  abstractPara.insertOoxml(\n    paragraphOoxml(\n      \"Abstract\",\n      \"Answers to the questions on trigonometry, using degrees to measure angles.\"\n    ),\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Split the Title paragraph and the Abstract paragraph into one run per\n# word, with each inter-word space as its own run (matching the \"one\n# word per <w:r>\" convention already used elsewhere in this document).\n#\n# A plain Range.Text assignment / InsertAfter rebuild would cause the\n# host to coalesce same-formatted adjacent runs back into a single run\n# on save, so instead we build the exact run-level OOXML for each\n# paragraph and swap it in with Range.InsertXML, which preserves run\n# boundaries verbatim.\n\nfunction Build-WordRunsXml($text) {\n    $words = $text -split ' '\n    $xml = \"\"\n    for ($i = 0; $i -lt $words.Length; $i++) {\n        if ($i -gt 0) {\n            $xml += '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>'\n        }\n        $word = $words[$i] -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n        $xml += '<w:r><w:t xml:space=\"preserve\">' + $word + '</w:t></w:r>'\n    }\n    return $xml\n}\n\nfunction Build-ParagraphOoxml($pStyle, $text) {\n    $runs = Build-WordRunsXml $text\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData>' + `\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n        '<w:body><w:p><w:pPr><w:pStyle w:val=\"' + $pStyle + '\"/></w:pPr>' + $runs + '</w:p></w:body></w:document>' + `\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$d = $word.ActiveDocument\n\n$titlePara = $null\n$abstractPara = $null\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Range.Style.NameLocal\n    if ($titlePara -eq $null -and $styleName -eq \"Title\") {\n        $titlePara = $p\n    }\n    if ($abstractPara -eq $null -and $styleName -eq \"Abstract\") {\n        $abstractPara = $p\n    }\n}\n\nif ($titlePara -ne $null) {\n    $xml = Build-ParagraphOoxml \"Title\" \"Answers: Trigonometry (degrees)\"\n    $null = $titlePara.Range.InsertXML($xml)\n}\n\nif ($abstractPara -ne $null) {\n    $xml = Build-ParagraphOoxml \"Abstract\" \"Answers to the questions on trigonometry, using degrees to measure angles.\"\n    $null = $abstractPara.Range.InsertXML($xml)\n}\n"}
